# "updated BNVP and AVMC with US versions"
#
# The "Cost Data" sheet's summary rows (88 and 96) were dividing the computed
# annual cost by 10 -- that adjustment is removed so the totals flow straight
# through as US-dollar annual values. The formatting for these three cells is
# also switched to match the style already used by the neighboring B87 total
# (percentage-style number format cell, no fill) instead of the old highlighted
# currency style.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cost Data")

# Remove the "/10" scaling that was applied to these rollup formulas.
$ws.Range("B88").Formula = "=B54"
$ws.Range("C88").Formula = "=B55"
$ws.Range("B96").Formula = "=B87"

# Match the cell formatting to B87 (same row group), which drops the old
# highlighted currency style these cells used to carry.
$ws.Range("B87").Copy()
$ws.Range("B88:C88").PasteSpecial(-4122)
$ws.Range("B96").PasteSpecial(-4122)
$excel.CutCopyMode = 0
